$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A22").Value = "Renewable"
$ws.Range("A21").Value = "Fossil Fuel"
$ws.Range("B21").Value = "#2A4845"
$ws.Range("B22").Value = "#41B496"

$ws.Range("G18").Select()
